$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the "favorites" and "what_does_it_mean" sheets, keeping "books".
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("favorites").Delete() | Out-Null
$wb.Worksheets.Item("what_does_it_mean").Delete() | Out-Null

$ws = $wb.Worksheets.Item("books")

# ---------------------------------------------------------------------------
# 2. Reset the header row formatting entirely (it was a bold custom-formatted
#    row) and then put the plain bordered look on C1:D1, reusing the
#    plain/bordered formatting that the regular data cells already use
#    (e.g. A16).
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).ClearFormats() | Out-Null

$ws.Range("A16").Copy() | Out-Null
$ws.Range("C1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Highlight column A (rows 1-15) in yellow and column B (rows 1-15) in
#    pink, reusing the existing fill styles already present on the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A1:A15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B5").Copy() | Out-Null
$ws.Range("B1:B15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Clear a couple of stray values.
# ---------------------------------------------------------------------------
$ws.Range("B1").ClearContents() | Out-Null
$ws.Range("D2").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 5. Fix up the duplicated/incorrectly styled ISBNs for rows 11 and 12 so
#    they match the formatting used by every other ISBN cell, and correct
#    the stray ISBN value left over in C17.
# ---------------------------------------------------------------------------
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C11:C12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C11").Value = 1285159454
$ws.Range("C12").Value = 1285159454
$ws.Range("C17").Value = 1285159454

# ---------------------------------------------------------------------------
# 6. Leave the cursor on C1.
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("C1").Select() | Out-Null
